$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: fill in the Laufzeit columns for test 4 / scenario "b" ---
# Write order follows the original author's column population sequence
# (G, H, I, K, J, L) so newly interned shared strings land in the same
# slots as the authored workbook.
$ws.Range("G15").Value = "01:02:40h"
$ws.Range("H15").Value = "01:38m"
$ws.Range("H15").NumberFormat = "h:mm"
$ws.Range("I15").Value = "00:03:26h"
$ws.Range("K15").Value = "00:59:36h"
$ws.Range("J15").Value = "01:00:55h"
$ws.Range("L15").Value = "2,71GB"

# --- Row 16: fill in the Laufzeit columns for test 4 / scenario "c" ---
# Write order: G, H, J, K, I, L
$ws.Range("G16").Value = "01:03:20h"
$ws.Range("H16").Value = "01:43m"
$ws.Range("H16").NumberFormat = "h:mm"
$ws.Range("J16").Value = "01:01:30h"
$ws.Range("K16").Value = "00:58:53h"
$ws.Range("I16").Value = "00:04:18h"
$ws.Range("I16").NumberFormat = "h:mm"
$ws.Range("L16").Value = "2,71GB"

# --- Row 18: new test 5, scenario "a" ---
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "a"
$ws.Range("C18").Value = 337
$ws.Range("D18").Value = 500
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 144
# Write order: G, H, I, K, J, L
$ws.Range("G18").Value = "02:23:56h"
$ws.Range("H18").Value = "03:18m"
$ws.Range("H18").NumberFormat = "h:mm"
$ws.Range("I18").Value = "00:05:14h"
$ws.Range("K18").Value = "02:18:35h"
$ws.Range("J18").Value = "02:20:36h"
$ws.Range("L18").Value = "2,71GB"

# --- Row 19: test 5, scenario "b" ---
$ws.Range("B19").Value = "b"
$ws.Range("C19").Value = 337
$ws.Range("D19").Value = 500
$ws.Range("E19").Value = 250
$ws.Range("F19").Value = 144
# Write order: G, H, I, J, K, L
$ws.Range("G19").Value = "02:26:06h"
$ws.Range("H19").Value = "01:46m"
$ws.Range("H19").NumberFormat = "h:mm"
$ws.Range("I19").Value = "00:02:53h"
$ws.Range("J19").Value = "02:24:12h"
$ws.Range("K19").Value = "02:23:04h"
$ws.Range("L19").Value = "2,71GB"
$ws.Range("M19").Value = 1

# --- Row 20: test 5, scenario "c" ---
$ws.Range("B20").Value = "c"
$ws.Range("C20").Value = 337
$ws.Range("D20").Value = 500
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 144

# --- Update selection to match the final authored state ---
$null = $ws.Range("L20").Select()
